# Apply the edit described by the diff:
# - Fill in X3 (PriceChange) and Y3 (UpDown) on the existing row 3
# - Append a new row 4 with a full set of values (no PriceChange/UpDown yet)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: add PriceChange / UpDown values
$ws.Range("X3").Value = -0.069999999999993179
$ws.Range("Y3").Value = "Down"

# Row 4: new data row (column A already carries a column-level date style,
# so it is inherited automatically; S/T need the percent style copied over
# from the row above since that style is only applied per-cell)
$ws.Range("A4").Value = 42641.890092592592
$ws.Range("B4").Value = -11
$ws.Range("C4").Value = "Sell"
$ws.Range("D4").Value = -16
$ws.Range("E4").Value = 21867
$ws.Range("F4").Value = 2342
$ws.Range("G4").Value = 55
$ws.Range("H4").Value = 45
$ws.Range("I4").Value = 39
$ws.Range("J4").Value = 60
$ws.Range("K4").Value = 10749
$ws.Range("L4").Value = 275
$ws.Range("M4").Value = 225
$ws.Range("N4").Value = 11
$ws.Range("O4").Value = 17
$ws.Range("P4").Value = "Bag"
$ws.Range("Q4").Value = 64.728146835133757
$ws.Range("R4").Value = -32.1

$ws.Range("S3").Copy()
$ws.Range("S4").PasteSpecial(-4122)
$ws.Range("S4").Value = -0.0755

$ws.Range("T3").Copy()
$ws.Range("T4").PasteSpecial(-4122)
$ws.Range("T4").Value = -0.0025

$ws.Range("U4").Value = 6.79
$ws.Range("V4").Value = 1.88
$ws.Range("W4").Value = 0
